# Harvard case classification: recalculated stats for the "_old" app
# columns (failure-type-ignored=no-disease-found), and inserted a new
# "average_doctor" column ahead of the legacy "average_doctor_old" one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: make room for the new average_doctor column by shifting
# the previous average_doctor values into average_doctor_old.
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4: updated recalculated stats
$ws.Range("E4").Value = 0.484
$ws.Range("F4").Value = 0.052
$ws.Range("G4").Value = 0.227
$ws.Range("N4").Value = 0.477
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("Q4").Value = 0.052
$ws.Range("R4").Value = 0.035
$ws.Range("S4").Value = 0.188
$ws.Range("W4").Value = 0.367
$ws.Range("AI4").Value = 0.395
$ws.Range("AJ4").Value = 0.094
$ws.Range("AK4").Value = 0.306
$ws.Range("AU4").Value = 0.245
$ws.Range("AW4").Value = 0.158
$ws.Range("BA4").Value = 2.041
$ws.Range("BB4").Value = 0.142
$ws.Range("BC4").Value = 0.377
$ws.Range("BG4").Value = 0.722
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.75
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.68
$ws.Range("BQ4").Value = 0.761

# Row 5: updated recalculated stats
$ws.Range("E5").Value = 0.602
$ws.Range("F5").Value = 0.06
$ws.Range("G5").Value = 0.244
$ws.Range("N5").Value = 0.721
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.276
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.014
$ws.Range("S5").Value = 0.12
$ws.Range("W5").Value = 0.331
$ws.Range("X5").Value = 0.099
$ws.Range("Y5").Value = 0.315
$ws.Range("AI5").Value = 0.388
$ws.Range("AJ5").Value = 0.083
$ws.Range("AK5").Value = 0.288
$ws.Range("AU5").Value = 0.456
$ws.Range("AV5").Value = 0.074
$ws.Range("AW5").Value = 0.272
$ws.Range("BA5").Value = 1.294
$ws.Range("BB5").Value = 0.073
$ws.Range("BC5").Value = 0.271
$ws.Range("BG5").Value = 0.379
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.521
$ws.Range("BN5").Value = 0.048
$ws.Range("BO5").Value = 0.219
$ws.Range("BP5").Value = 0.431
$ws.Range("BQ5").Value = 0.455

# Row 6: updated recalculated stats
$ws.Range("E6").Value = 0.537
$ws.Range("N6").Value = 0.574
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.348
$ws.Range("AI6").Value = 0.391
$ws.Range("AU6").Value = 0.319
$ws.Range("BA6").Value = 1.578
$ws.Range("BG6").Value = 0.497
$ws.Range("BM6").Value = 0.615
$ws.Range("BP6").Value = 0.526
$ws.Range("BQ6").Value = 0.5669999999999999

# Row 7: updated recalculated stats
$ws.Range("E7").Value = 0.574
$ws.Range("N7").Value = 0.654
$ws.Range("Q7").Value = 0.037
$ws.Range("W7").Value = 0.338
$ws.Range("AI7").Value = 0.389
$ws.Range("AU7").Value = 0.389
$ws.Range("BA7").Value = 1.394
$ws.Range("BG7").Value = 0.419
$ws.Range("BM7").Value = 0.555
$ws.Range("BP7").Value = 0.465
$ws.Range("BQ7").Value = 0.494

# Row 8: updated recalculated stats
$ws.Range("E8").Value = 0.703
$ws.Range("F8").Value = 0.07199999999999999
$ws.Range("G8").Value = 0.268
$ws.Range("N8").Value = 0.802
$ws.Range("O8").Value = 0.06
$ws.Range("P8").Value = 0.244
$ws.Range("Q8").Value = 0.038
$ws.Range("W8").Value = 0.395
$ws.Range("X8").Value = 0.118
$ws.Range("Y8").Value = 0.343
$ws.Range("AI8").Value = 0.46
$ws.Range("AJ8").Value = 0.134
$ws.Range("AK8").Value = 0.366
$ws.Range("AU8").Value = 0.403
$ws.Range("AV8").Value = 0.081
$ws.Range("AW8").Value = 0.285
$ws.Range("BA8").Value = 1.754
$ws.Range("BB8").Value = 0.108
$ws.Range("BC8").Value = 0.329
$ws.Range("BG8").Value = 0.5580000000000001
$ws.Range("BH8").Value = 0.109
$ws.Range("BI8").Value = 0.331
$ws.Range("BM8").Value = 0.671
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.248
$ws.Range("BP8").Value = 0.585
$ws.Range("BQ8").Value = 0.621

# Row 9: updated recalculated stats
$ws.Range("E9").Value = 0.667
$ws.Range("F9").Value = 0.222
$ws.Range("G9").Value = 0.471
$ws.Range("N9").Value = 0.738
$ws.Range("O9").Value = 0.193
$ws.Range("P9").Value = 0.44
$ws.Range("W9").Value = 0.286
$ws.Range("X9").Value = 0.204
$ws.Range("Y9").Value = 0.452
$ws.Range("AI9").Value = 0.405
$ws.Range("AJ9").Value = 0.241
$ws.Range("AK9").Value = 0.491
$ws.Range("BA9").Value = 1.714
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.595
$ws.Range("BH9").Value = 0.241
$ws.Range("BI9").Value = 0.491
$ws.Range("BM9").Value = 0.643
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.571
$ws.Range("BQ9").Value = 0.612

# Row 10: updated recalculated stats
$ws.Range("E10").Value = 0.8100000000000001
$ws.Range("F10").Value = 0.154
$ws.Range("G10").Value = 0.393
$ws.Range("N10").Value = 0.929
$ws.Range("O10").Value = 0.066
$ws.Range("P10").Value = 0.258
$ws.Range("W10").Value = 0.5
$ws.Range("X10").Value = 0.25
$ws.Range("Y10").Value = 0.5
$ws.Range("AI10").Value = 0.5
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.405
$ws.Range("AV10").Value = 0.241
$ws.Range("AW10").Value = 0.491
$ws.Range("BA10").Value = 2.167
$ws.Range("BB10").Value = 0.214
$ws.Range("BC10").Value = 0.462
$ws.Range("BG10").Value = 0.667
$ws.Range("BH10").Value = 0.222
$ws.Range("BI10").Value = 0.471
$ws.Range("BM10").Value = 0.8100000000000001
$ws.Range("BN10").Value = 0.154
$ws.Range("BO10").Value = 0.393
$ws.Range("BP10").Value = 0.722
$ws.Range("BQ10").Value = 0.758

# Row 11: updated recalculated stats
$ws.Range("E11").Value = 0.857
$ws.Range("F11").Value = 0.122
$ws.Range("G11").Value = 0.35
$ws.Range("N11").Value = 0.929
$ws.Range("O11").Value = 0.066
$ws.Range("P11").Value = 0.258
$ws.Range("W11").Value = 0.5
$ws.Range("X11").Value = 0.25
$ws.Range("Y11").Value = 0.5
$ws.Range("AI11").Value = 0.571
$ws.Range("AJ11").Value = 0.245
$ws.Range("AK11").Value = 0.495
$ws.Range("AU11").Value = 0.571
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.167
$ws.Range("BB11").Value = 0.214
$ws.Range("BC11").Value = 0.462
$ws.Range("BG11").Value = 0.667
$ws.Range("BH11").Value = 0.222
$ws.Range("BI11").Value = 0.471
$ws.Range("BM11").Value = 0.8100000000000001
$ws.Range("BN11").Value = 0.154
$ws.Range("BO11").Value = 0.393
$ws.Range("BP11").Value = 0.722
$ws.Range("BQ11").Value = 0.764

# Row 12: updated recalculated stats
$ws.Range("E12").Value = 1.417
$ws.Range("F12").Value = 0.854
$ws.Range("G12").Value = 0.924
$ws.Range("N12").Value = 1.256
$ws.Range("O12").Value = 0.293
$ws.Range("P12").Value = 0.542
$ws.Range("W12").Value = 1.524
$ws.Range("X12").Value = 0.44
$ws.Range("Y12").Value = 0.663
$ws.Range("AI12").Value = 1.625
$ws.Range("AJ12").Value = 1.484
$ws.Range("AK12").Value = 1.218
$ws.Range("AU12").Value = 2.846
$ws.Range("AV12").Value = 3.361
$ws.Range("AW12").Value = 1.833
$ws.Range("BA12").Value = 3.792
$ws.Range("BB12").Value = 0.449
$ws.Range("BC12").Value = 0.67
$ws.Range("BG12").Value = 1.143
$ws.Range("BH12").Value = 0.194
$ws.Range("BI12").Value = 0.44
$ws.Range("BM12").Value = 1.235
$ws.Range("BN12").Value = 0.239
$ws.Range("BO12").Value = 0.489
$ws.Range("BP12").Value = 1.264
$ws.Range("BQ12").Value = 1.249

# Row 13: updated recalculated stats
$ws.Range("E13").Value = 1.379
$ws.Range("F13").Value = 0.27
$ws.Range("G13").Value = 0.519
$ws.Range("N13").Value = 1.742
$ws.Range("O13").Value = 0.476
$ws.Range("P13").Value = 0.6899999999999999
$ws.Range("W13").Value = 0.972
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.438
$ws.Range("AI13").Value = 1.144
$ws.Range("AJ13").Value = 0.309
$ws.Range("AK13").Value = 0.556
$ws.Range("AU13").Value = 2.039
$ws.Range("AV13").Value = 0.339
$ws.Range("AW13").Value = 0.582
$ws.Range("BA13").Value = 2.175
$ws.Range("BB13").Value = 0.281
$ws.Range("BC13").Value = 0.53
$ws.Range("BG13").Value = 0.543
$ws.Range("BH13").Value = 0.051
$ws.Range("BI13").Value = 0.226
$ws.Range("BM13").Value = 0.783
$ws.Range("BN13").Value = 0.166
$ws.Range("BO13").Value = 0.407
$ws.Range("BP13").Value = 0.725
$ws.Range("BQ13").Value = 0.664
